$wb = $excel.ActiveWorkbook

# Sheet ALC, row 28 (Leve Item ID 27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 50001656
$ws.Range("J28").Value = 83334950
$ws.Range("L28").Value = 83334950
$ws.Range("N28").Value = -83335920

# Sheet ALC, row 62 (Leve Item ID 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2890.05
$ws.Range("I62").Value = 2513.0667
$ws.Range("K62").Value = 2513.0667
$ws.Range("M62").Value = -1889.0667

# Sheet ALC, row 64 (Leve Item ID 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5397.2856
$ws.Range("I64").Value = 5136.1875
$ws.Range("J64").Value = 6232.8
$ws.Range("K64").Value = 5136.1875
$ws.Range("L64").Value = 6232.8
$ws.Range("M64").Value = -4888.1875
$ws.Range("N64").Value = -6728.8

# Sheet ALC, row 65 (Leve Item ID 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2890.05
$ws.Range("I65").Value = 2513.0667
$ws.Range("K65").Value = 12565.3335
$ws.Range("M65").Value = -9445.333499999999

# Sheet ALC, row 67 (Leve Item ID 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5397.2856
$ws.Range("I67").Value = 5136.1875
$ws.Range("J67").Value = 6232.8
$ws.Range("K67").Value = 5136.1875
$ws.Range("L67").Value = 6232.8
$ws.Range("M67").Value = -4278.1875
$ws.Range("N67").Value = -7948.8

# Sheet ALC, row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4622.4
$ws.Range("I76").Value = 3263.8
$ws.Range("J76").Value = 5981
$ws.Range("K76").Value = 3263.8
$ws.Range("L76").Value = 5981
$ws.Range("M76").Value = -2948.8
$ws.Range("N76").Value = -6611

# Sheet ALC, row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4622.4
$ws.Range("I79").Value = 3263.8
$ws.Range("J79").Value = 5981
$ws.Range("K79").Value = 3263.8
$ws.Range("L79").Value = 5981
$ws.Range("M79").Value = -2171.8
$ws.Range("N79").Value = -8165

# Sheet ALC, row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7109.5
$ws.Range("I113").Value = 5549.3335
$ws.Range("J113").Value = 9449.75
$ws.Range("K113").Value = 5549.3335
$ws.Range("L113").Value = 9449.75
$ws.Range("M113").Value = -2295.3335
$ws.Range("N113").Value = -15957.75

# Sheet ALC, row 116 (Leve Item ID 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 15664.833
$ws.Range("I116").Value = 13497
$ws.Range("K116").Value = 13497
$ws.Range("M116").Value = -10055

# Sheet ARM, row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4770.9
$ws.Range("I110").Value = 5054.1113
$ws.Range("K110").Value = 5054.1113
$ws.Range("M110").Value = -3009.1113

# Sheet BSM, row 62 (Leve Item ID 10586)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 199500
$ws.Range("J62").Value = 199500
$ws.Range("L62").Value = 199500
$ws.Range("N62").Value = -200872

# Sheet BSM, row 65 (Leve Item ID 10586)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 199500
$ws.Range("J65").Value = 199500
$ws.Range("L65").Value = 598500
$ws.Range("N65").Value = -605364

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1098.5
$ws.Range("J86").Value = 1098.5
$ws.Range("L86").Value = 1098.5
$ws.Range("N86").Value = -3344.5

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1098.5
$ws.Range("J89").Value = 1098.5
$ws.Range("L89").Value = 5492.5
$ws.Range("N89").Value = -16724.5

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9154.549999999999
$ws.Range("I31").Value = 3227.2778
$ws.Range("J31").Value = 62500
$ws.Range("K31").Value = 3227.2778
$ws.Range("L31").Value = 62500
$ws.Range("M31").Value = -2932.2778
$ws.Range("N31").Value = -63090

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9154.549999999999
$ws.Range("I34").Value = 3227.2778
$ws.Range("J34").Value = 62500
$ws.Range("K34").Value = 3227.2778
$ws.Range("L34").Value = 62500
$ws.Range("M34").Value = -3025.2778
$ws.Range("N34").Value = -62904

# Sheet CRP, row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4390
$ws.Range("I58").Value = 1200
$ws.Range("K58").Value = 1200
$ws.Range("M58").Value = -997

# Sheet CRP, row 122 (Leve Item ID 36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2388.138
$ws.Range("I122").Value = 2083.1904
$ws.Range("J122").Value = 3188.625
$ws.Range("K122").Value = 6249.5712
$ws.Range("L122").Value = 9565.875
$ws.Range("M122").Value = -3799.5712
$ws.Range("N122").Value = -14465.875

# Sheet CRP, row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4390
$ws.Range("I136").Value = 1200
$ws.Range("K136").Value = 3600
$ws.Range("M136").Value = -1050

# Sheet CUL, row 4 (Leve Item ID 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20806944
$ws.Range("J4").Value = 381602.1
$ws.Range("L4").Value = 1144806.3
$ws.Range("N4").Value = -1145030.3

# Sheet CUL, row 12 (Leve Item ID 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 457.08334
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 497.72726
$ws.Range("K12").Value = 30
$ws.Range("L12").Value = 1493.18178
$ws.Range("M12").Value = 143
$ws.Range("N12").Value = -1839.18178

# Sheet CUL, row 14 (Leve Item ID 12886)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 5082.909
$ws.Range("I14").Value = 5082.909
$ws.Range("K14").Value = 15248.727
$ws.Range("M14").Value = -15075.727

# Sheet CUL, row 23 (Leve Item ID 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 139.33333
$ws.Range("I23").Value = 83.333336
$ws.Range("J23").Value = 195.33333
$ws.Range("K23").Value = 250.000008
$ws.Range("L23").Value = 585.99999
$ws.Range("M23").Value = -15.00000800000001
$ws.Range("N23").Value = -1055.99999

# Sheet CUL, row 62 (Leve Item ID 12845)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5236
$ws.Range("I62").Value = 1755.6666
$ws.Range("J62").Value = 10456.5
$ws.Range("K62").Value = 5266.9998
$ws.Range("L62").Value = 31369.5
$ws.Range("M62").Value = -4580.9998
$ws.Range("N62").Value = -32741.5

# Sheet CUL, row 65 (Leve Item ID 12845)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 5236
$ws.Range("I65").Value = 1755.6666
$ws.Range("J65").Value = 10456.5
$ws.Range("K65").Value = 15800.9994
$ws.Range("L65").Value = 94108.5
$ws.Range("M65").Value = -12368.9994
$ws.Range("N65").Value = -100972.5

# Sheet CUL, row 86 (Leve Item ID 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 6249.5557
$ws.Range("J86").Value = 828
$ws.Range("L86").Value = 2484
$ws.Range("N86").Value = -4856

# Sheet CUL, row 89 (Leve Item ID 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 6249.5557
$ws.Range("J89").Value = 828
$ws.Range("L89").Value = 7452
$ws.Range("N89").Value = -19308

# Sheet CUL, row 138 (Leve Item ID 44105)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2743.889
$ws.Range("I138").Value = 2661.2856
$ws.Range("K138").Value = 7983.8568
$ws.Range("M138").Value = -2843.8568

# Sheet CUL, row 139 (Leve Item ID 44102)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2747.7778
$ws.Range("J139").Value = 3033
$ws.Range("L139").Value = 9099
$ws.Range("N139").Value = -19379

# Sheet CUL, row 140 (Leve Item ID 44097)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1631.5416
$ws.Range("J140").Value = 1802
$ws.Range("L140").Value = 5406
$ws.Range("N140").Value = -15766

# Sheet GSM, row 70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6149.5713
$ws.Range("I70").Value = 4426.467
$ws.Range("K70").Value = 4426.467
$ws.Range("M70").Value = -4156.467

# Sheet GSM, row 73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6149.5713
$ws.Range("I73").Value = 4426.467
$ws.Range("K73").Value = 4426.467
$ws.Range("M73").Value = -3490.467

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4405.5
$ws.Range("J122").Value = 10749.5
$ws.Range("L122").Value = 32248.5
$ws.Range("N122").Value = -37148.5

# Sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4322.8203
$ws.Range("I40").Value = 3608.0303
$ws.Range("K40").Value = 3608.0303
$ws.Range("M40").Value = -3472.0303

# Sheet LTW, row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4755.3687
$ws.Range("I122").Value = 4081.9395
$ws.Range("J122").Value = 9200
$ws.Range("K122").Value = 12245.8185
$ws.Range("L122").Value = 27600
$ws.Range("M122").Value = -9795.818499999999
$ws.Range("N122").Value = -32500

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9070.714
$ws.Range("I132").Value = 3372.5
$ws.Range("K132").Value = 10117.5
$ws.Range("M132").Value = -7587.5

# Sheet WVR, row 107 (Leve Item ID 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1177765.6
$ws.Range("I107").Value = 2223624.5
$ws.Range("J107").Value = 1174.375
$ws.Range("K107").Value = 6670873.5
$ws.Range("L107").Value = 3523.125
$ws.Range("M107").Value = -6668953.5
$ws.Range("N107").Value = -7363.125

# Sheet WVR, row 122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2003.1621
$ws.Range("J122").Value = 2486.5454
$ws.Range("L122").Value = 7459.6362
$ws.Range("N122").Value = -12359.6362
